# OCEPROJECT-4912 - added load 'click' test for R4R results
#
# Adds four new "click" scenario rows (all / tools / areas / filtered) to the
# R4RResultsLoad sheet, with two new columns (ActionStatus, Filters) that
# describe the simulated click action + querystring filters used by the load
# test, and renames the existing baseline row's label from "R4R Results" to
# "R4R Results (all)". Also switches the active sheet/tab from R4RHomeLoad to
# R4RResultsLoad.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("R4RResultsLoad")

# --- Header row (row 1): add ActionStatus / Filters columns -----------------
$ws.Cells.Item(1, 3).Value = "ActionStatus"
$ws.Cells.Item(1, 4).Value = "Filters"

# Copy the bold header formatting from A1 onto the two new header cells so
# they share the existing "header" cell style instead of creating a new one.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C1:D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 2: existing baseline row, renamed + extended ------------------------
$ws.Cells.Item(2, 1).Value = "/research/resources/search"
$ws.Cells.Item(2, 2).Value = "R4R Results (all)"
$ws.Cells.Item(2, 3).Value = "r4r_results|view|none|ra=0;tt=0;rt=0;tst=0|1|"
$ws.Cells.Item(2, 4).Value = "none"

# --- Row 3: tools filter click scenario --------------------------------------
$ws.Cells.Item(3, 1).Value = "/research/resources/search?from=0&toolTypes=analysis_tools"
$ws.Cells.Item(3, 2).Value = "R4R Results (tools)"
$ws.Cells.Item(3, 3).Value = "r4r_results|view|none|ra=0;tt=1;rt=0;tst=0|1|"
$ws.Cells.Item(3, 4).Value = "analysis_tools"

# --- Row 4: research-areas filter click scenario -----------------------------
$ws.Cells.Item(4, 1).Value = "/research/resources/search?from=0&researchAreas=cancer_omics"
$ws.Cells.Item(4, 2).Value = "R4R Results (areas)"
$ws.Cells.Item(4, 3).Value = "r4r_results|view|none|ra=1;tt=0;rt=0;tst=0|1|"
$ws.Cells.Item(4, 4).Value = "cancer_omics"

# --- Row 5: combined/filtered (paged) click scenario -------------------------
$ws.Cells.Item(5, 1).Value = "/research/resources/search?from=20&toolSubtypes=modeling&toolSubtypes=r_software&toolTypes=analysis_tools"
$ws.Cells.Item(5, 2).Value = "R4R Results (filtered)"
$ws.Cells.Item(5, 3).Value = "r4r_results|view|none|ra=0;tt=1;rt=0;tst=2|2|"
$ws.Cells.Item(5, 4).Value = "modeling|r_software|analysis_tools"

# --- Column widths (approximate best-fit for the new/changed content) -------
$ws.Columns.Item(1).ColumnWidth = 106
$ws.Columns.Item(2).ColumnWidth = 19.1666666666667
$ws.Columns.Item(3).ColumnWidth = 42.1666666666667
$ws.Columns.Item(4).ColumnWidth = 33.5

# --- Activate R4RResultsLoad as the selected tab, with A6 selected ----------
$ws.Activate() | Out-Null
$ws.Range("A6").Select() | Out-Null

Write-Host "R4R Results load-test rows added"
